# Updates cryptocurrency price/volume data in the worksheet.
# Applies Price (column D) and Volume(1h) (column E) updates for rows 2-51.
# Numeric-looking text values are entered with a leading apostrophe so that
# Excel keeps them as text (matching the original inline-string cell type)
# rather than silently converting/reformatting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.797.43'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '3.178.59'
$ws.Range("E3").Value = '  -4.76%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'" + '571.25'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = "'" + '171.96'
$ws.Range("E6").Value = '  -3.06%  '
$ws.Range("D7").Value = "'" + '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'" + '0.601'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").Value = '3.180.66'
$ws.Range("E9").Value = '  -4.51%  '
$ws.Range("E10").Value = '  -3.06%  '
$ws.Range("D11").Value = "'" + '6.57'
$ws.Range("E11").Value = '  -4.40%  '
$ws.Range("D12").Value = "'" + '0.392'
$ws.Range("E12").Value = '  -4.73%  '
$ws.Range("D13").Value = '3.720.22'
$ws.Range("E13").Value = '  -4.78%  '
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").Value = "'" + '27.40'
$ws.Range("E15").Value = '  -4.60%  '
$ws.Range("D16").Value = '65.666.81'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = "'" + '0.0000164'
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").Value = '3.171.76'
$ws.Range("E18").Value = '  -4.86%  '
$ws.Range("D19").Value = "'" + '5.71'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = "'" + '12.90'
$ws.Range("E20").Value = '  -3.66%  '
$ws.Range("D21").Value = "'" + '359.94'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").Value = "'" + '7.28'
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = "'" + '69.13'
$ws.Range("E24").Value = '  -2.91%  '
$ws.Range("D25").Value = "'" + '0.495'
$ws.Range("E25").Value = '  -4.78%  '
$ws.Range("D26").Value = '3.305.53'
$ws.Range("E26").Value = '  -5.24%  '
$ws.Range("D27").Value = "'" + '0.0000115'
$ws.Range("E27").Value = '  -5.53%  '
$ws.Range("D28").Value = "'" + '9.85'
$ws.Range("E28").Value = '  +2.77%  '
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = "'" + '1.00'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("D33").Value = "'" + '5.38'
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("D34").Value = "'" + '21.97'
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("E35").Value = '  -1.51%  '
$ws.Range("D36").Value = "'" + '6.62'
$ws.Range("E36").Value = '  -3.47%  '
$ws.Range("D37").Value = "'" + '159.94'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("D39").Value = "'" + '0.838'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").Value = "'" + '1.80'
$ws.Range("E40").Value = '  +2.83%  '
$ws.Range("D41").Value = "'" + '26.41'
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").Value = "'" + '2.52'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = '2.647.69'
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").Value = "'" + '6.12'
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("D45").Value = "'" + '4.19'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = "'" + '0.0659'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").Value = "'" + '329.39'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = "'" + '24.13'
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("D50").Value = "'" + '0.0274'
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("E51").Value = '  -0.80%  '
